# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (interested-attendee count) figures in column F
# for sheets "展览" (Worksheets index 1) and "全部类型" (Worksheets index 4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws1.Range("F7").Value  = 7744
$ws1.Range("F11").Value = 6386
$ws1.Range("F12").Value = 3311
$ws1.Range("F24").Value = 3725
$ws1.Range("F29").Value = 1382
$ws1.Range("F32").Value = 2682
$ws1.Range("F33").Value = 1667
$ws1.Range("F37").Value = 3465
$ws1.Range("F38").Value = 237
$ws1.Range("F42").Value = 510
$ws1.Range("F43").Value = 1343
$ws1.Range("F46").Value = 614

$ws4 = $wb.Worksheets.Item(4)   # 全部类型
$ws4.Range("F10").Value = 7744
$ws4.Range("F13").Value = 6386
$ws4.Range("F14").Value = 3311
$ws4.Range("F24").Value = 3725
$ws4.Range("F31").Value = 1383
$ws4.Range("F34").Value = 2682
$ws4.Range("F35").Value = 1667
$ws4.Range("F40").Value = 3465
$ws4.Range("F41").Value = 237
$ws4.Range("F45").Value = 510
$ws4.Range("F46").Value = 1343
$ws4.Range("F49").Value = 614
